# Fixed bug with removing commas within quotes
# Corrects several truncated words in the "entries" sheet's entry_vocal_part
# column (F) that were accidentally cut short (commas inside quoted strings
# previously caused the rest of the text to be dropped).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("entries")

$ws.Range("F4").Value  = [char]0x201D + "Bass" + [char]0x201D + ", A"
$ws.Range("F11").Value = [char]0x201D + "Cant" + [char]0x201D + ", Am"
$ws.Range("F12").Value = "att. " + [char]0x201C + "Milgrove" + [char]0x201D + ", 3 voices"
$ws.Range("F13").Value = [char]0x201C + "Rippon" + [char]0x2019 + "s Coll.n" + [char]0x201D + ", 3 voices"
$ws.Range("F15").Value = [char]0x201C + "Dr. Addington" + [char]0x2019 + "s Colln." + [char]0x201D + ", 3 voices"
$ws.Range("F16").Value = "att. " + [char]0x201C + "Dr Arne" + [char]0x201D + ", 3 voices"
$ws.Range("F17").Value = "att. " + [char]0x201C + "I. Smith" + [char]0x201D + ", 3 voices"
$ws.Range("F18").Value = "att. " + [char]0x201C + "Leach" + [char]0x201D + ", 4 voices"
$ws.Range("F24").Value = [char]0x201C + "Bass" + [char]0x201D + ", D"
